$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (new weekly data inserted at the top of the dataset)
$ws.Range("D2").Value = 44496
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 2800
$ws.Range("T2").Value = 10

# Row 3 (previously row 2's data)
$ws.Range("D3").Value = 44483
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 35
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 5

# Row 4 (previously row 3's data)
$ws.Range("D4").Value = 44488
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("R4").Value = "La Ligua"
$ws.Range("S4").Value = 2400
$ws.Range("T4").Value = 5

# Row 5 (previously row 4's data)
$ws.Range("D5").Value = 44166
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "La Ligua"
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 18

# Row 6 (new row, contains what used to be row 5's data)
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44466
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100104
$ws.Range("H6").Value = "Frutos de pepita"
$ws.Range("I6").Value = 100104004
$ws.Range("J6").Value = "Níspero"
$ws.Range("K6").Value = "Californiana(o)"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 11000
$ws.Range("Q6").Value = "$/bandeja 5 kilos"
$ws.Range("R6").Value = "La Ligua"
$ws.Range("S6").Value = 2200
$ws.Range("T6").Value = 5

Write-Host "Done"
